$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/IF1 should carry the same formatting as the other
# header cells (bold, centered, bordered) -- copy H1's format instead of
# re-building it by hand so the existing style (s="1") is reused verbatim.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I is a constant-1 column, J duplicates H (IP) for every data row.
for ($r = 2; $r -le 29; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
